# Auto-generated edit script: updates crypto price/volume table cells
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.234.22"
$ws.Range("E2").Value = "  +0.42%  "
$ws.Range("D3").Value = "1.589.01"
$ws.Range("E3").Value = "  +1.05%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.37"
$ws.Range("E5").Value = "  +1.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.501"
$ws.Range("E6").Value = "  +0.34%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  +0.63%  "
$ws.Range("E9").Value = "  -0.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.32"
$ws.Range("E10").Value = "  -1.21%  "
$ws.Range("E11").Value = "  +0.77%  "
$ws.Range("D12").Value = "1.813.44"
$ws.Range("E12").Value = "  +1.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.01"
$ws.Range("E13").Value = "  -0.86%  "
$ws.Range("D14").Value = "1.553.13"
$ws.Range("E14").Value = "  -1.36%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.520"
$ws.Range("E15").Value = "  +1.19%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.27"
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("D17").Value = "26.249.48"
$ws.Range("E17").Value = "  +0.52%  "
$ws.Range("E18").Value = "  +0.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.40"
$ws.Range("E19").Value = "  +1.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "213.39"
$ws.Range("E20").Value = "  +2.71%  "
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.28"
$ws.Range("E22").Value = "  +0.99%  "
$ws.Range("B23").Value = "Toncoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.17"
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("B24").Value = "Avalanche"
$ws.Range("C24").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.03"
$ws.Range("E24").Value = "  +2.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.71"
$ws.Range("E25").Value = "  +0.18%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.05"
$ws.Range("E27").Value = "  +1.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.111"
$ws.Range("E28").Value = "  -0.40%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.18"
$ws.Range("E29").Value = "  -0.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0496"
$ws.Range("E30").Value = "  -1.79%  "
$ws.Range("E31").Value = "  +1.59%  "
$ws.Range("E32").Value = "  -0.15%  "
$ws.Range("D33").Value = "1.344.98"
$ws.Range("E33").Value = "  +5.35%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.93"
$ws.Range("E34").Value = "  -1.39%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.45"
$ws.Range("E35").Value = "  +0.26%  "
$ws.Range("E36").Value = "  -0.22%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.586"
$ws.Range("E37").Value = "  -3.97%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0167"
$ws.Range("E38").Value = "  +1.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.821"
$ws.Range("E39").Value = "  +1.38%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.74"
$ws.Range("E40").Value = "  +3.47%  "
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.978"
$ws.Range("E42").Value = "  -10.04%  "
$ws.Range("E43").Value = "  +0.72%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.766"
$ws.Range("E44").Value = "  +0.65%  "
$ws.Range("D45").Value = "1.724.88"
$ws.Range("E45").Value = "  +1.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "61.43"
$ws.Range("E46").Value = "  -1.56%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "85.94"
$ws.Range("E47").Value = "  -3.31%  "
$ws.Range("D48").Value = "0.0₆0104"
$ws.Range("E48").Value = "  +4.82%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.47"
$ws.Range("E49").Value = "  -2.41%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0977"
$ws.Range("E50").Value = "  -2.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0501"
$ws.Range("E51").Value = "  -0.92%  "
